{"js": "// Convert every \"simple field\" (<w:fldSimple w:instr=\"...\"/>) that still\n// uses the shorthand OOXML syntax into the equivalent, Word-canonical\n// \"complex field\" run sequence:\n//   <w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>\n//   <w:r><w:instrText>INSTR</w:instrText></w:r>\n//   <w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>\n//   <w:r><w:fldChar w:fldCharType=\"end\"/></w:r>\n// Everything else in the owning paragraph (its own formatting, any\n// bookmarks, etc.) is left untouched.\n\nfunction xmlEscapeText(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nfunction xmlUnescapeAttr(s) {\n  return s\n    .replace(/&quot;/g, '\"')\n    .replace(/&apos;/g, \"'\")\n    .replace(/&lt;/g, \"<\")\n    .replace(/&gt;/g, \">\")\n    .replace(/&amp;/g, \"&\");\n}\n\nfunction fieldRunsXml(instr) {\n  const code = xmlEscapeText(xmlUnescapeAttr(instr));\n  return (\n    '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n    \"<w:r><w:instrText>\" + code + \"</w:instrText></w:r>\" +\n    '<w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n    '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>'\n  );\n}\n\n// Replace <w:fldSimple w:instr=\"...\">...</w:fldSimple> / self-closing form\n// with the run-based expansion; returns null if nothing changed.\nfunction expandFldSimple(paragraphXml) {\n  let changed = false;\n  let out = paragraphXml.replace(\n    /<w:fldSimple\\s+w:instr=\"([^\"]*)\"\\s*\\/>/g,\n    (whole, instr) => {\n      changed = true;\n      return fieldRunsXml(instr);\n    }\n  );\n  out = out.replace(\n    /<w:fldSimple\\s+w:instr=\"([^\"]*)\"\\s*>([\\s\\S]*?)<\\/w:fldSimple>/g,\n    (whole, instr) => {\n      changed = true;\n      return fieldRunsXml(instr);\n    }\n  );\n  return changed ? out : null;\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find every paragraph that hosts at least one field (fldSimple elements\n// surface as Field objects through the Office.js object model too).\nconst candidates = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const contentRange = paragraph.getRange(\"Content\");\n  const fields = contentRange.fields;\n  fields.load(\"items\");\n  candidates.push({ paragraph, fields });\n}\nawait context.sync();\n\nconst targets = candidates.filter((c) => c.fields.items.length > 0);\n\n// Pull the exact OOXML (including the paragraph's own attributes and any\n// sibling content such as bookmarks) for every target paragraph.\nconst withOoxml = targets.map((t) => ({\n  paragraph: t.paragraph,\n  ooxml: t.paragraph.getRange(\"Whole\").getOoxml(),\n}));\nawait context.sync();\n\nfor (const entry of withOoxml) {\n  const full = entry.ooxml.value;\n  const bodyMatch = full.match(/<w:body>([\\s\\S]*)<\\/w:body>/);\n  if (!bodyMatch) {\n    continue;\n  }\n  // The OOXML roundtrip always returns the target paragraph followed by a\n  // trailing empty paragraph used as an anchor; only the first <w:p>...</w:p>\n  // is the paragraph we actually want to rewrite.\n  const pMatch = bodyMatch[1].match(/<w:p\\b[\\s\\S]*?<\\/w:p>|<w:p\\b[^>]*\\/>/);\n  if (!pMatch) {\n    continue;\n  }\n  let paragraphXml = pMatch[0];\n  // Strip the synthetic w14:paraId / w14:textId attributes the OOXML\n  // roundtrip adds; they are not present in the original document.\n  paragraphXml = paragraphXml\n    .replace(/\\s+w14:paraId=\"[^\"]*\"/, \"\")\n    .replace(/\\s+w14:textId=\"[^\"]*\"/, \"\");\n\n  const expanded = expandFldSimple(paragraphXml);\n  if (!expanded) {\n    continue;\n  }\n  // Re-supply the paragraph's own attributes (rsids, etc.) on the wrapper\n  // <w:p> of the replacement package - even though we are inserting into the\n  // paragraph's *content* range (so the paragraph mark itself is not part of\n  // the replaced range), this engine mirrors the supplied <w:p> attributes\n  // back onto the live paragraph rather than leaving them untouched.\n  const packageXml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" + expanded + \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\";\n\n  entry.paragraph.getRange(\"Content\").insertOoxml(packageXml, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Convert every \"simple field\" (<w:fldSimple w:instr=\"...\"/>) that still\n# uses the shorthand OOXML syntax into the equivalent, Word-canonical\n# \"complex field\" run sequence:\n#   <w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>\n#   <w:r><w:instrText>INSTR</w:instrText></w:r>\n#   <w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>\n#   <w:r><w:fldChar w:fldCharType=\"end\"/></w:r>\n# Everything else in the owning paragraph (its own formatting, any\n# bookmarks, etc.) is left untouched.\n\n$d = $word.ActiveDocument\n\n$fldSimpleSelfClosing = '<w:fldSimple\\s+w:instr=\"([^\"]*)\"\\s*/>'\n$fldSimpleWithBody = '<w:fldSimple\\s+w:instr=\"([^\"]*)\"\\s*>.*?</w:fldSimple>'\n$fieldRunsReplacement = '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r><w:r><w:instrText>$1</w:instrText></w:r><w:r><w:fldChar w:fldCharType=\"separate\"/></w:r><w:r><w:fldChar w:fldCharType=\"end\"/></w:r>'\n\n$pkgHeader = \"<?xml version=`\"1.0`\"?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>\"\n$pkgFooter = \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n\n    if ($r.Fields.Count -lt 1) {\n        continue\n    }\n\n    $xml = $r.WordOpenXML\n    if (-not ($xml -match '<w:body>([\\s\\S]*)</w:body>')) {\n        continue\n    }\n    $bodyInner = $matches[1]\n\n    if (-not ($bodyInner -match '(<w:p\\b[\\s\\S]*?</w:p>|<w:p\\b[^>]*/>)')) {\n        continue\n    }\n    $paraXml = $matches[1]\n\n    # Strip the synthetic w14:paraId / w14:textId attributes the OOXML\n    # roundtrip adds; they are not present in the original document.\n    $paraXml = $paraXml -replace '\\s+w14:paraId=\"[^\"]*\"', ''\n    $paraXml = $paraXml -replace '\\s+w14:textId=\"[^\"]*\"', ''\n\n    if ($paraXml -notmatch '<w:fldSimple\\s+w:instr=\"[^\"]*\"\\s*/?>') {\n        continue\n    }\n\n    $newParaXml = $paraXml -replace $fldSimpleSelfClosing, $fieldRunsReplacement\n    $newParaXml = $newParaXml -replace $fldSimpleWithBody, $fieldRunsReplacement\n\n    $ooxml = $pkgHeader + $newParaXml + $pkgFooter\n    $r.InsertXML($ooxml)\n}\n"}
